# Turn deref into a factory
# - add a new "no_deref" test row (A4/B4) to the TestSheet
# - document the new case with a comment on B4 (mirrors the existing
#   B2/B3 "key/fallback/parser" comment convention)
# - leave the cursor where it was left at save time

$excel.UserName = "Microsoft Office User"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: A4 = "no_deref", B4 = 20
$ws.Range("A4").Value = "no_deref"
$ws.Range("B4").Value = 20

# Comment documenting the deref-factory "no" case for B4
$commentText = "Microsoft Office User:`n{{--`nkey: <<A4>>`nfallback: 0`nparser: int`nderef: no`n--}}"
$comment = $ws.Range("B4").AddComment($commentText)

# Match the rich-text look of the sibling comments (bold header line,
# larger body font) as closely as the host lets us.
$comment.Shape.TextFrame.Characters.Font.Name = "Calibri"
$comment.Shape.TextFrame.Characters.Font.Size = 18
$comment.Shape.TextFrame.Characters(1, 23).Font.Bold = $true

# Restore the selection recorded in the saved workbook
$ws.Range("Q11").Select()
